{"js": "// Vietnamese translation pass for \"Email 3 [TEMPLATE] Partner email \u2013 list of travel documents\"\n// Each entry replaces the FIRST remaining occurrence of `find` with `replace`.\n// Several source phrases repeat (e.g. \" on \", \" or \"), so the list below\n// lists each occurrence separately, in document order, and we consume one\n// match at a time so later duplicate occurrences keep their own translation.\nconst replacements = [\n  [\"Thank you for registering for \", \"C\u1ea3m \u01a1n b\u1ea1n \u0111\u00e3 \u0111\u0103ng k\u00fd tham d\u1ef1 s\u1ef1 ki\u1ec7n \"],\n  [\"Hi \", \"Xin ch\u00e0o \"],\n  [\"We are excited for you to join us at \", \"Ch\u00fang t\u00f4i r\u1ea5t vui khi b\u1ea1n c\u00f3 th\u1ec3 tham d\u1ef1 v\u1edbi ch\u00fang t\u00f4i t\u1ea1i s\u1ef1 ki\u1ec7n \"],\n  [\n    \"To confirm your registration, we would require you and one guest of your choice to provide us with:\",\n    \"\u0110\u1ec3 ti\u1ebfn h\u00e0nh \u0111\u0103ng k\u00fd cho b\u1ea1n, ch\u00fang t\u00f4i c\u1ea7n b\u1ea1n v\u00e0 m\u1ed9t kh\u00e1ch m\u1eddi c\u1ee7a b\u1ea1n cung c\u1ea5p cho ch\u00fang t\u00f4i:\",\n  ],\n  [\"A scanned copy of your international passports\", \"B\u1ea3n scan h\u1ed9 chi\u1ebfu qu\u1ed1c t\u1ebf c\u1ee7a b\u1ea1n\"],\n  [\"Covid-19 vaccination certificates\", \"Gi\u1ea5y ch\u1ee9ng nh\u1eadn ti\u00eam ph\u00f2ng Covid-19\"],\n  [\"Send my details\", \"G\u1eedi th\u00f4ng tin c\u1ee7a t\u00f4i\"],\n  [\n    \"Your country manager will be in touch to confirm your booking or request any other relevant details. \",\n    \"Gi\u00e1m \u0111\u1ed1c ph\u1ee5 tr\u00e1ch t\u1ea1i qu\u1ed1c gia c\u1ee7a b\u1ea1n s\u1ebd li\u00ean l\u1ea1c \u0111\u1ec3 x\u00e1c nh\u1eadn ho\u1eb7c h\u1ecfi th\u00eam c\u00e1c th\u00f4ng tin li\u00ean quan kh\u00e1c n\u1ebfu c\u1ea7n. \",\n  ],\n  [\"Our event package offers you and your guest: \", \"G\u00f3i s\u1ef1 ki\u1ec7n ch\u00fang t\u00f4i cung c\u1ea5p \u0111\u1ebfn b\u1ea1n v\u00e0 kh\u00e1ch m\u1eddi c\u1ee7a b\u1ea1n bao g\u1ed3m: \"],\n  [\"Flight tickets \", \"V\u00e9 m\u00e1y bay \"],\n  [\"Travel insurance \", \"B\u1ea3o hi\u1ec3m du l\u1ecbch \"],\n  [\"Airport \u2013 Hotel \u2013 Airport transfer \", \"\u0110\u01b0a \u0111\u00f3n s\u00e2n bay \u2013 kh\u00e1ch s\u1ea1n \"],\n  [\n    \"One hotel room for you and your guest / Two hotel rooms for you and your guest\",\n    \"M\u1ed9t ho\u1eb7c hai ph\u00f2ng kh\u00e1ch s\u1ea1n cho b\u1ea1n v\u00e0 kh\u00e1ch m\u1eddi c\u1ee7a b\u1ea1n\",\n  ],\n  [\"Check-in\", \"Nh\u1eadn ph\u00f2ng\"],\n  [\" on \", \" v\u00e0o ng\u00e0y \"],\n  [\"Check-out\", \"Tr\u1ea3 ph\u00f2ng\"],\n  [\" on \", \" v\u00e0o ng\u00e0y \"],\n  [\"Meals (Breakfast, lunch, and dinner)\", \"C\u00e1c b\u1eefa \u0103n (B\u1eefa s\u00e1ng, b\u1eefa tr\u01b0a v\u00e0 b\u1eefa t\u1ed1i)\"],\n  [\n    \"We will send you a confirmation letter before your departure date with the event agenda and information about your flights, transportation, and accommodation. \",\n    \"Ch\u00fang t\u00f4i s\u1ebd g\u1eedi th\u01b0 x\u00e1c nh\u1eadn \u0111\u1ebfn b\u1ea1n tr\u01b0\u1edbc ng\u00e0y kh\u1edfi h\u00e0nh v\u1edbi c\u00e1c th\u00f4ng tin chi ti\u1ebft v\u1ec1 ch\u01b0\u01a1ng tr\u00ecnh s\u1ef1 ki\u1ec7n, chuy\u1ebfn bay, ph\u01b0\u01a1ng ti\u1ec7n di chuy\u1ec3n v\u00e0 ch\u1ed7 \u1edf c\u1ee7a b\u1ea1n. \",\n  ],\n  [\"If you have any questions, please contact us via \", \"N\u1ebfu b\u1ea1n c\u1ea7n h\u1ed7 tr\u1ee3, vui l\u00f2ng li\u00ean h\u1ec7 v\u1edbi ch\u00fang t\u00f4i qua \"],\n  [\" or \", \" ho\u1eb7c \"],\n  [\n    \"If you have any questions, please contact your country manager, \",\n    \"N\u1ebfu b\u1ea1n c\u00f3 b\u1ea5t k\u1ef3 th\u1eafc m\u1eafc n\u00e0o, vui l\u00f2ng li\u00ean h\u1ec7 v\u1edbi gi\u00e1m \u0111\u1ed1c ph\u1ee5 tr\u00e1ch qu\u1ed1c gia c\u1ee7a b\u1ea1n \",\n  ],\n  [\", at \", \", qua email \"],\n  [\" or \", \" ho\u1eb7c s\u1ed1 \"],\n  [\"We look forward to seeing you soon.\", \"Ch\u00fang t\u00f4i r\u1ea5t mong \u0111\u01b0\u1ee3c g\u1eb7p b\u1ea1n.\"],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + find);\n  }\n  results.items[0].insertText(replace, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Vietnamese translation pass for \"Email 3 [TEMPLATE] Partner email \u2013 list of travel documents\"\n# Walks the document top-to-bottom, replacing one occurrence of each English\n# phrase with its Vietnamese translation using Find/Replace (wdReplaceOne).\n# A couple of English phrases repeat verbatim (\" on \", \" or \") but need\n# different Vietnamese text at each spot, so every replacement is scoped to\n# a Range that starts right after the end of the previous replacement -\n# this guarantees we always hit the next untouched occurrence, in order.\n\n$d = $word.ActiveDocument\n\nfunction ReplaceNext($findText, $replaceText, $searchStart) {\n  $r = $d.Range($searchStart, $d.Content.End)\n  $r.Find.ClearFormatting()\n  $r.Find.Replacement.ClearFormatting()\n  $ok = $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1)\n  if (-not $ok) {\n    throw \"ReplaceNext: could not find '$findText'\"\n  }\n  return $r.End\n}\n\n$pos = 0\n$pos = ReplaceNext \"Thank you for registering for \" \"C\u1ea3m \u01a1n b\u1ea1n \u0111\u00e3 \u0111\u0103ng k\u00fd tham d\u1ef1 s\u1ef1 ki\u1ec7n \" $pos\n$pos = ReplaceNext \"Hi \" \"Xin ch\u00e0o \" $pos\n$pos = ReplaceNext \"We are excited for you to join us at \" \"Ch\u00fang t\u00f4i r\u1ea5t vui khi b\u1ea1n c\u00f3 th\u1ec3 tham d\u1ef1 v\u1edbi ch\u00fang t\u00f4i t\u1ea1i s\u1ef1 ki\u1ec7n \" $pos\n$pos = ReplaceNext \"To confirm your registration, we would require you and one guest of your choice to provide us with:\" \"\u0110\u1ec3 ti\u1ebfn h\u00e0nh \u0111\u0103ng k\u00fd cho b\u1ea1n, ch\u00fang t\u00f4i c\u1ea7n b\u1ea1n v\u00e0 m\u1ed9t kh\u00e1ch m\u1eddi c\u1ee7a b\u1ea1n cung c\u1ea5p cho ch\u00fang t\u00f4i:\" $pos\n$pos = ReplaceNext \"A scanned copy of your international passports\" \"B\u1ea3n scan h\u1ed9 chi\u1ebfu qu\u1ed1c t\u1ebf c\u1ee7a b\u1ea1n\" $pos\n$pos = ReplaceNext \"Covid-19 vaccination certificates\" \"Gi\u1ea5y ch\u1ee9ng nh\u1eadn ti\u00eam ph\u00f2ng Covid-19\" $pos\n$pos = ReplaceNext \"Send my details\" \"G\u1eedi th\u00f4ng tin c\u1ee7a t\u00f4i\" $pos\n$pos = ReplaceNext \"Your country manager will be in touch to confirm your booking or request any other relevant details. \" \"Gi\u00e1m \u0111\u1ed1c ph\u1ee5 tr\u00e1ch t\u1ea1i qu\u1ed1c gia c\u1ee7a b\u1ea1n s\u1ebd li\u00ean l\u1ea1c \u0111\u1ec3 x\u00e1c nh\u1eadn ho\u1eb7c h\u1ecfi th\u00eam c\u00e1c th\u00f4ng tin li\u00ean quan kh\u00e1c n\u1ebfu c\u1ea7n. \" $pos\n$pos = ReplaceNext \"Our event package offers you and your guest: \" \"G\u00f3i s\u1ef1 ki\u1ec7n ch\u00fang t\u00f4i cung c\u1ea5p \u0111\u1ebfn b\u1ea1n v\u00e0 kh\u00e1ch m\u1eddi c\u1ee7a b\u1ea1n bao g\u1ed3m: \" $pos\n$pos = ReplaceNext \"Flight tickets \" \"V\u00e9 m\u00e1y bay \" $pos\n$pos = ReplaceNext \"Travel insurance \" \"B\u1ea3o hi\u1ec3m du l\u1ecbch \" $pos\n$pos = ReplaceNext \"Airport \u2013 Hotel \u2013 Airport transfer \" \"\u0110\u01b0a \u0111\u00f3n s\u00e2n bay \u2013 kh\u00e1ch s\u1ea1n \" $pos\n$pos = ReplaceNext \"One hotel room for you and your guest / Two hotel rooms for you and your guest\" \"M\u1ed9t ho\u1eb7c hai ph\u00f2ng kh\u00e1ch s\u1ea1n cho b\u1ea1n v\u00e0 kh\u00e1ch m\u1eddi c\u1ee7a b\u1ea1n\" $pos\n$pos = ReplaceNext \"Check-in\" \"Nh\u1eadn ph\u00f2ng\" $pos\n$pos = ReplaceNext \" on \" \" v\u00e0o ng\u00e0y \" $pos\n$pos = ReplaceNext \"Check-out\" \"Tr\u1ea3 ph\u00f2ng\" $pos\n$pos = ReplaceNext \" on \" \" v\u00e0o ng\u00e0y \" $pos\n$pos = ReplaceNext \"Meals (Breakfast, lunch, and dinner)\" \"C\u00e1c b\u1eefa \u0103n (B\u1eefa s\u00e1ng, b\u1eefa tr\u01b0a v\u00e0 b\u1eefa t\u1ed1i)\" $pos\n$pos = ReplaceNext \"We will send you a confirmation letter before your departure date with the event agenda and information about your flights, transportation, and accommodation. \" \"Ch\u00fang t\u00f4i s\u1ebd g\u1eedi th\u01b0 x\u00e1c nh\u1eadn \u0111\u1ebfn b\u1ea1n tr\u01b0\u1edbc ng\u00e0y kh\u1edfi h\u00e0nh v\u1edbi c\u00e1c th\u00f4ng tin chi ti\u1ebft v\u1ec1 ch\u01b0\u01a1ng tr\u00ecnh s\u1ef1 ki\u1ec7n, chuy\u1ebfn bay, ph\u01b0\u01a1ng ti\u1ec7n di chuy\u1ec3n v\u00e0 ch\u1ed7 \u1edf c\u1ee7a b\u1ea1n. \" $pos\n$pos = ReplaceNext \"If you have any questions, please contact us via \" \"N\u1ebfu b\u1ea1n c\u1ea7n h\u1ed7 tr\u1ee3, vui l\u00f2ng li\u00ean h\u1ec7 v\u1edbi ch\u00fang t\u00f4i qua \" $pos\n$pos = ReplaceNext \" or \" \" ho\u1eb7c \" $pos\n$pos = ReplaceNext \"If you have any questions, please contact your country manager, \" \"N\u1ebfu b\u1ea1n c\u00f3 b\u1ea5t k\u1ef3 th\u1eafc m\u1eafc n\u00e0o, vui l\u00f2ng li\u00ean h\u1ec7 v\u1edbi gi\u00e1m \u0111\u1ed1c ph\u1ee5 tr\u00e1ch qu\u1ed1c gia c\u1ee7a b\u1ea1n \" $pos\n$pos = ReplaceNext \", at \" \", qua email \" $pos\n$pos = ReplaceNext \" or \" \" ho\u1eb7c s\u1ed1 \" $pos\n$pos = ReplaceNext \"We look forward to seeing you soon.\" \"Ch\u00fang t\u00f4i r\u1ea5t mong \u0111\u01b0\u1ee3c g\u1eb7p b\u1ea1n.\" $pos\n"}
